$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2043.75
$ws.Range("J51").Value = 2043.75
$ws.Range("L51").Value = 2043.75
$ws.Range("N51").Value = -3011.75
# row 51: A Bile Business / Shark Oil

$ws.Range("H62").Value = 2856.2222
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 3034.3333
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 3034.3333
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -4282.3333
# row 62: The Mustache Suits Him / Enchanted Mythrite Ink

$ws.Range("H64").Value = 2938.7222
$ws.Range("I64").Value = 2799.6
$ws.Range("J64").Value = 2992.2307
$ws.Range("K64").Value = 2799.6
$ws.Range("L64").Value = 2992.2307
$ws.Range("M64").Value = -2551.6
$ws.Range("N64").Value = -3488.2307
# row 64: Forged from the Void / Void Glue

$ws.Range("H65").Value = 2856.2222
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 3034.3333
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 15171.6665
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -21411.6665
# row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink

$ws.Range("H67").Value = 2938.7222
$ws.Range("I67").Value = 2799.6
$ws.Range("J67").Value = 2992.2307
$ws.Range("K67").Value = 2799.6
$ws.Range("L67").Value = 2992.2307
$ws.Range("M67").Value = -1941.6
$ws.Range("N67").Value = -4708.2307
# row 67: Dodging the Draft (L) / Void Glue

$ws.Range("H135").Value = 1313
$ws.Range("I135").Value = 629.25
$ws.Range("J135").Value = 2407
$ws.Range("K135").Value = 5663.25
$ws.Range("L135").Value = 21663
$ws.Range("M135").Value = -3128.25
$ws.Range("N135").Value = -26733
# row 135: For Tired Minds / Grade 1 Gemsap of Intelligence

$ws.Range("H137").Value = 1718.48
$ws.Range("I137").Value = 2415.1428
$ws.Range("J137").Value = 1447.5555
$ws.Range("K137").Value = 7245.428400000001
$ws.Range("L137").Value = 4342.666499999999
$ws.Range("M137").Value = -4695.428400000001
$ws.Range("N137").Value = -9442.666499999999
# row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1327176.8
$ws.Range("I32").Value = 1522723.1
$ws.Range("K32").Value = 1522723.1
$ws.Range("M32").Value = -1522436.1
# row 32: Ingot We Trust / Steel Ingot

$ws.Range("H61").Value = 9807123
$ws.Range("I61").Value = 17546488
$ws.Range("K61").Value = 17546488
$ws.Range("M61").Value = -17546276
# row 61: Dealing with the Tough Stuff / Cobalt Ingot

$ws.Range("H74").Value = 13891832
$ws.Range("I74").Value = 1440.2
$ws.Range("J74").Value = 23813540
$ws.Range("K74").Value = 1440.2
$ws.Range("L74").Value = 23813540
$ws.Range("M74").Value = -566.2
$ws.Range("N74").Value = -23815288
# row 74: As the Bolt Flies / Titanium Nugget

$ws.Range("H77").Value = 13891832
$ws.Range("I77").Value = 1440.2
$ws.Range("J77").Value = 23813540
$ws.Range("K77").Value = 7201
$ws.Range("L77").Value = 119067700
$ws.Range("M77").Value = -2833
$ws.Range("N77").Value = -119076436
# row 77: Heavy Metal Banned (L) / Titanium Nugget

$ws.Range("H97").Value = 934.03845
$ws.Range("I97").Value = 719.1667
$ws.Range("J97").Value = 1417.5
$ws.Range("K97").Value = 719.1667
$ws.Range("L97").Value = 1417.5
$ws.Range("M97").Value = -223.1667
$ws.Range("N97").Value = -2409.5
# row 97: Ore for Me / High Steel Ingot

$ws.Range("H136").Value = 9807123
$ws.Range("I136").Value = 17546488
$ws.Range("K136").Value = 52639464
$ws.Range("M136").Value = -52636914
# row 136: Metal with Mettle / Cobalt Tungsten Ingot

$ws.Range("H138").Value = 53200
$ws.Range("J138").Value = 53200
$ws.Range("L138").Value = 53200
$ws.Range("N138").Value = -63480
# row 138: Don't Ask about the Rivets / Titanium Gold Helm of Casting

$ws.Range("H139").Value = 69919.09
$ws.Range("J139").Value = 67911
$ws.Range("L139").Value = 67911
$ws.Range("N139").Value = -78191
# row 139: Backing up My Words / Titanium Gold Thornplate of Fending

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2012.4546
$ws.Range("I99").Value = 1768
$ws.Range("J99").Value = 2104.125
$ws.Range("K99").Value = 1768
$ws.Range("L99").Value = 2104.125
$ws.Range("M99").Value = -270
$ws.Range("N99").Value = -5100.125
# row 99: O Pine / Pine Lumber

$ws.Range("H122").Value = 1531.9333
$ws.Range("I122").Value = 1187.091
$ws.Range("J122").Value = 1731.579
$ws.Range("K122").Value = 3561.273
$ws.Range("L122").Value = 5194.737
$ws.Range("M122").Value = -1111.273
$ws.Range("N122").Value = -10094.737
# row 122: Timber of Tenkonto / Horse Chestnut Lumber

$ws.Range("H126").Value = 2012.4546
$ws.Range("I126").Value = 1768
$ws.Range("J126").Value = 2104.125
$ws.Range("K126").Value = 5304
$ws.Range("L126").Value = 6312.375
$ws.Range("M126").Value = -2834
$ws.Range("N126").Value = -11252.375
# row 126: A Better Conductor / Red Pine Lumber

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4157.143
$ws.Range("I80").Value = 4350
$ws.Range("J80").Value = 4080
$ws.Range("K80").Value = 13050
$ws.Range("L80").Value = 12240
$ws.Range("M80").Value = -12114
$ws.Range("N80").Value = -14112
# row 80: Saucy for a Suitor / Hollandaise Sauce

$ws.Range("H83").Value = 4157.143
$ws.Range("I83").Value = 4350
$ws.Range("J83").Value = 4080
$ws.Range("K83").Value = 39150
$ws.Range("L83").Value = 36720
$ws.Range("M83").Value = -34470
$ws.Range("N83").Value = -46080
# row 83: Saved by the Sauce (L) / Hollandaise Sauce

$ws.Range("H137").Value = 7505.92
$ws.Range("I137").Value = 17802.715
$ws.Range("J137").Value = 3501.611
$ws.Range("K137").Value = 53408.145
$ws.Range("L137").Value = 10504.833
$ws.Range("M137").Value = -48308.145
$ws.Range("N137").Value = -20704.833
# row 137: Creative Chocolate / Gateau au Chocolat

$ws.Range("H140").Value = 1358.5333
$ws.Range("I140").Value = 1134.1428
$ws.Range("K140").Value = 3402.4284
$ws.Range("M140").Value = 1777.5716
# row 140: Sweet, Sweet Bean Juice / Mesquite Juice

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2267.25
$ws.Range("I126").Value = 2376.8572
$ws.Range("K126").Value = 7130.571599999999
$ws.Range("M126").Value = -4660.571599999999
# row 126: Gold Rush Order / Phrygian Gold Ingot

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 59012.832
$ws.Range("J28").Value = 59012.832
$ws.Range("L28").Value = 59012.832
$ws.Range("N28").Value = -59708.832
# row 28: Doublet Jeopardy / Cotton Doublet Vest of Gathering

$ws.Range("H62").Value = 121025.37
$ws.Range("I62").Value = 201555.8
$ws.Range("J62").Value = 53916.668
$ws.Range("K62").Value = 201555.8
$ws.Range("L62").Value = 53916.668
$ws.Range("M62").Value = -200931.8
$ws.Range("N62").Value = -55164.668
# row 62: Pride Up in Smoke / Rainbow Cloth

$ws.Range("H65").Value = 121025.37
$ws.Range("I65").Value = 201555.8
$ws.Range("J65").Value = 53916.668
$ws.Range("K65").Value = 1007779
$ws.Range("L65").Value = 269583.34
$ws.Range("M65").Value = -1004659
$ws.Range("N65").Value = -275823.34
# row 65: Desperate for Diversionaries (L) / Rainbow Cloth

$ws.Range("H132").Value = 6947350
$ws.Range("I132").Value = 2949.48
$ws.Range("K132").Value = 8848.440000000001
$ws.Range("M132").Value = -6318.440000000001
# row 132: Comfy Cabins / Snow Cotton Cloth
